$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row and data rows 1-5 with the revised requirement content ---
$ws.Range("A1").Value = "Orden de navegación"
$ws.Range("B1").Value = "Atributo en pantalla"
$ws.Range("C1").Value = "Tipo"
$ws.Range("D1").Value = "Descripción"
$ws.Range("E1").Value = "Aspectos a tener en cuenta"
$ws.Range("F1").Value = "Tooltip"
$ws.Range("G1").Value = "Validaciones en el ingreso"
$ws.Range("H1").Value = "Mensaje de error"
$ws.Range("I1").Value = "Parámetros de mensajes"
$ws.Range("J1").Value = "Obligatorio"
$ws.Range("K1").Value = "Permite modificarse"
$ws.Range("L1").Value = "Valor por defecto"
$ws.Range("M1").Value = "Tipo de campo"
$ws.Range("N1").Value = "Número máximo de caracteres"
$ws.Range("O1").Value = "Número mínimo de caracteres"
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Correo electrónico"
$ws.Range("C2").Value = "Input text"
$ws.Range("D2").Value = "Campo para el ingreso del correo electrónico del usuario"
$ws.Range("E2").Value = "El usuario ya debe estar registrado en el sistema."
$ws.Range("F2").Value = "N/A"
$ws.Range("G2").Value = "Usuario existente"
$ws.Range("H2").Value = "1. El correo electrónico es obligatorio.                     2. Por favor ingresa un correo electrónico válido. "
$ws.Range("I2").Value = "N/A"
$ws.Range("J2").Value = "Si"
$ws.Range("K2").Value = "Si"
$ws.Range("L2").Value = "Ninguno"
$ws.Range("M2").Value = "Alfanumérico"
$ws.Range("N2").Value = "N/A"
$ws.Range("O2").Value = "N/A"
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Contraseña"
$ws.Range("C3").Value = "Input password"
$ws.Range("D3").Value = "Campo para el ingreso de contraseña correspondinete al correo ingresado."
$ws.Range("E3").Value = "La contraseña debe coincidir con la asociada al correo electrónico ingresado anteriormente"
$ws.Range("F3").Value = "N/A"
$ws.Range("G3").Value = "Puede recibir caracteres alfanuméricos y especiales"
$ws.Range("H3").Value = "1. La contraseña es obligatoria.                     2. Por favor ingrese más de 9 caracteres."
$ws.Range("I3").Value = "N/A"
$ws.Range("J3").Value = "Si"
$ws.Range("K3").Value = "Si"
$ws.Range("L3").Value = "Ninguno"
$ws.Range("M3").Value = "Alfanumérico"
$ws.Range("N3").Value = "N/A"
$ws.Range("O3").Value = 10
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Ingresar"
$ws.Range("C4").Value = "Button"
$ws.Range("D4").Value = "Al ser presionado el sistema valida que los datos ingresados en los campos de Correo electrónico y contraseña sean correctos; si ese es el caso, el usuario ingresa a la apliación respectiva al tipo de usuario. En caso contrario se debe mostrar el error con el mensaje que informe el tipo de este."
$ws.Range("E4").Value = "Los campos de usuario y contraseña deben ser correctos"
$ws.Range("F4").Value = "Iniciar sesión"
$ws.Range("G4").Value = "El campo de correo electrónico y contraseña deben contener texto"
$ws.Range("H4").Value = "1. Correo electrónico o contraseña incorrectos.  2. Correo electrónico no encontrado en el sistema."
$ws.Range("I4").Value = "N/A"
$ws.Range("J4").Value = "N/A"
$ws.Range("K4").Value = "N/A"
$ws.Range("L4").Value = "N/A"
$ws.Range("M4").Value = "N/A"
$ws.Range("N4").Value = "N/A"
$ws.Range("O4").Value = "N/A"
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Registrarse"
$ws.Range("C5").Value = "Button"
$ws.Range("D5").Value = "Al ser presionado mostrará una ventana en la cual el usuario podrá realizar el registro en el sistema"
$ws.Range("E5").Value = "Ninguno"
$ws.Range("F5").Value = "Registrarse en el sistema"
$ws.Range("G5").Value = "No hay validaciones"
$ws.Range("H5").Value = "N/A"
$ws.Range("I5").Value = "N/A"
$ws.Range("J5").Value = "N/A"
$ws.Range("K5").Value = "N/A"
$ws.Range("L5").Value = "N/A"
$ws.Range("M5").Value = "N/A"
$ws.Range("N5").Value = "N/A"
$ws.Range("O5").Value = "N/A"

# --- Remove the obsolete numeric-digit column (old column P) ---
$ws.Columns("P").Delete()

# --- Remove the duplicated "Registrarme" row (old row 6) ---
$ws.Rows(6).Delete()

# --- Restore the active cell selection ---
$ws.Range("K4").Select()
